# Update the "Förändrad" (changed) date column (C) for all data rows
# from serial date 45179 (2023-09-10) to 45180 (2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 422
}

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45180
